# Update specific cost/value cells in column C per the commit's data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value  = 40.98
$ws.Range("C8").Value  = 24.87
$ws.Range("C10").Value = 109.27
$ws.Range("C13").Value = 71.08
$ws.Range("C14").Value = 101.93
$ws.Range("C15").Value = 250.32
$ws.Range("C16").Value = 84.58
$ws.Range("C18").Value = 102.73
$ws.Range("C20").Value = 115.17
$ws.Range("C21").Value = 74.3
$ws.Range("C22").Value = 104.67
$ws.Range("C23").Value = 96.47
